$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# --- Neg_Change sheet: replace data rows 2-14 (13 rows, same row count) ---
$ws1.Range("A2").Value = "TECHM"
$ws1.Range("B2").Value = 1445.3
$ws1.Range("C2").Value = 1464.7
$ws1.Range("D2").Value = 1430.7
$ws1.Range("E2").Value = 1460
$ws1.Range("F2").Value = 2818268
$ws1.Range("G2").Value = 5777445
$ws1.Range("H2").Value = -0.512194750447646
$ws1.Range("I2").Value = "TECHM"
$ws1.Range("A3").Value = "MAXHEALTH"
$ws1.Range("B3").Value = 1166
$ws1.Range("C3").Value = 1172
$ws1.Range("D3").Value = 1147.9
$ws1.Range("E3").Value = 1168.8
$ws1.Range("F3").Value = 2784798
$ws1.Range("G3").Value = 5964536
$ws1.Range("H3").Value = -0.5331073531956216
$ws1.Range("I3").Value = "MAXHEALTH"
$ws1.Range("A4").Value = "WIPRO"
$ws1.Range("B4").Value = 247.5
$ws1.Range("C4").Value = 247.8
$ws1.Range("D4").Value = 245.71
$ws1.Range("E4").Value = 246.38
$ws1.Range("F4").Value = 6705862
$ws1.Range("G4").Value = 13298452
$ws1.Range("H4").Value = -0.4957411584446069
$ws1.Range("I4").Value = "WIPRO"
$ws1.Range("A5").Value = "APOLLOHOSP"
$ws1.Range("B5").Value = 7481
$ws1.Range("C5").Value = 7491
$ws1.Range("D5").Value = 7415
$ws1.Range("E5").Value = 7440
$ws1.Range("F5").Value = 228411
$ws1.Range("G5").Value = 454714
$ws1.Range("H5").Value = -0.4976820594923402
$ws1.Range("I5").Value = "APOLLOHOSP"
$ws1.Range("A6").Value = "TITAN"
$ws1.Range("B6").Value = 3950
$ws1.Range("C6").Value = 3956
$ws1.Range("D6").Value = 3895
$ws1.Range("E6").Value = 3900
$ws1.Range("F6").Value = 762391
$ws1.Range("G6").Value = 1578123
$ws1.Range("H6").Value = -0.5169001402298807
$ws1.Range("I6").Value = "TITAN"
$ws1.Range("A7").Value = "HCLTECH"
$ws1.Range("B7").Value = 1662.6
$ws1.Range("C7").Value = 1668.4
$ws1.Range("D7").Value = 1642.6
$ws1.Range("E7").Value = 1644.4
$ws1.Range("F7").Value = 2724035
$ws1.Range("G7").Value = 6342672
$ws1.Range("H7").Value = -0.5705224864221262
$ws1.Range("I7").Value = "HCLTECH"
$ws1.Range("A8").Value = "LTIM"
$ws1.Range("B8").Value = 6000
$ws1.Range("C8").Value = 6091
$ws1.Range("D8").Value = 5959
$ws1.Range("E8").Value = 6030
$ws1.Range("F8").Value = 398886
$ws1.Range("G8").Value = 868072
$ws1.Range("H8").Value = -0.5404920329189283
$ws1.Range("I8").Value = "LTIM"
$ws1.Range("A9").Value = "MOTHERSON"
$ws1.Range("B9").Value = 112.75
$ws1.Range("C9").Value = 113.39
$ws1.Range("D9").Value = 111.39
$ws1.Range("E9").Value = 112
$ws1.Range("F9").Value = 15140498
$ws1.Range("G9").Value = 33900889
$ws1.Range("H9").Value = -0.55338935211994
$ws1.Range("I9").Value = "MOTHERSON"
$ws1.Range("A10").Value = "SBICARD"
$ws1.Range("B10").Value = 863.7
$ws1.Range("C10").Value = 877.35
$ws1.Range("D10").Value = 863.7
$ws1.Range("E10").Value = 873.75
$ws1.Range("F10").Value = 390374
$ws1.Range("G10").Value = 862548
$ws1.Range("H10").Value = -0.5474176509597147
$ws1.Range("I10").Value = "SBICARD"
$ws1.Range("A11").Value = "OBEROIRLTY"
$ws1.Range("B11").Value = 1715
$ws1.Range("C11").Value = 1724.9
$ws1.Range("D11").Value = 1698
$ws1.Range("E11").Value = 1709.9
$ws1.Range("F11").Value = 154835
$ws1.Range("G11").Value = 334415
$ws1.Range("H11").Value = -0.5369974432965029
$ws1.Range("I11").Value = "OBEROIRLTY"
$ws1.Range("A12").Value = "SONACOMS"
$ws1.Range("B12").Value = 509.95
$ws1.Range("C12").Value = 515.5
$ws1.Range("D12").Value = 502.25
$ws1.Range("E12").Value = 506.5
$ws1.Range("F12").Value = 2647276
$ws1.Range("G12").Value = 5842354
$ws1.Range("H12").Value = -0.5468819588816426
$ws1.Range("I12").Value = "SONACOMS"
$ws1.Range("A13").Value = "PPLPHARMA"
$ws1.Range("B13").Value = 190
$ws1.Range("C13").Value = 192
$ws1.Range("D13").Value = 189.11
$ws1.Range("E13").Value = 190.15
$ws1.Range("F13").Value = 1601706
$ws1.Range("G13").Value = 3241928
$ws1.Range("H13").Value = -0.5059402923198788
$ws1.Range("I13").Value = "PPLPHARMA"
$ws1.Range("A14").Value = "KAYNES"
$ws1.Range("B14").Value = 6010
$ws1.Range("C14").Value = 6034.5
$ws1.Range("D14").Value = 5951.5
$ws1.Range("E14").Value = 5965.5
$ws1.Range("F14").Value = 335154
$ws1.Range("G14").Value = 702111
$ws1.Range("H14").Value = -0.5226481282874076
$ws1.Range("I14").Value = "KAYNES"

# --- Pos_Change sheet: replace data rows 2-10, and rows 11-15 are new/extended ---
$ws2.Range("A2").Value = "LT"
$ws2.Range("B2").Value = 4035
$ws2.Range("C2").Value = 4048
$ws2.Range("D2").Value = 4007.7
$ws2.Range("E2").Value = 4035.1
$ws2.Range("F2").Value = 1463574
$ws2.Range("G2").Value = 923032
$ws2.Range("H2").Value = 0.5856156666291092
$ws2.Range("I2").Value = "LT"
$ws2.Range("A3").Value = "INDHOTEL"
$ws2.Range("B3").Value = 722.4
$ws2.Range("C3").Value = 741
$ws2.Range("D3").Value = 719.1
$ws2.Range("E3").Value = 732.2
$ws2.Range("F3").Value = 4075549
$ws2.Range("G3").Value = 2796978
$ws2.Range("H3").Value = 0.457125869420496
$ws2.Range("I3").Value = "INDHOTEL"
$ws2.Range("A4").Value = "AMBUJACEM"
$ws2.Range("B4").Value = 562
$ws2.Range("C4").Value = 562
$ws2.Range("D4").Value = 555
$ws2.Range("E4").Value = 556.4
$ws2.Range("F4").Value = 1248817
$ws2.Range("G4").Value = 847533
$ws2.Range("H4").Value = 0.4734730093105519
$ws2.Range("I4").Value = "AMBUJACEM"
$ws2.Range("A5").Value = "VEDL"
$ws2.Range("B5").Value = 516
$ws2.Range("C5").Value = 519
$ws2.Range("D5").Value = 509.15
$ws2.Range("E5").Value = 510
$ws2.Range("F5").Value = 5736356
$ws2.Range("G5").Value = 4009545
$ws2.Range("H5").Value = 0.4306750516579811
$ws2.Range("I5").Value = "VEDL"
$ws2.Range("A6").Value = "LICI"
$ws2.Range("B6").Value = 915
$ws2.Range("C6").Value = 919.95
$ws2.Range("D6").Value = 906.05
$ws2.Range("E6").Value = 907.3
$ws2.Range("F6").Value = 1301140
$ws2.Range("G6").Value = 923863
$ws2.Range("H6").Value = 0.4083689897744579
$ws2.Range("I6").Value = "LICI"
$ws2.Range("A7").Value = "MPHASIS"
$ws2.Range("B7").Value = 2736.8
$ws2.Range("C7").Value = 2758.9
$ws2.Range("D7").Value = 2702.9
$ws2.Range("E7").Value = 2745
$ws2.Range("F7").Value = 1216126
$ws2.Range("G7").Value = 814207
$ws2.Range("H7").Value = 0.4936324546460544
$ws2.Range("I7").Value = "MPHASIS"
$ws2.Range("A8").Value = "ASHOKLEY"
$ws2.Range("B8").Value = 146.08
$ws2.Range("C8").Value = 146.85
$ws2.Range("D8").Value = 144.86
$ws2.Range("E8").Value = 146.85
$ws2.Range("F8").Value = 8625173
$ws2.Range("G8").Value = 5986604
$ws2.Range("H8").Value = 0.4407455378708864
$ws2.Range("I8").Value = "ASHOKLEY"
$ws2.Range("A9").Value = "PETRONET"
$ws2.Range("B9").Value = 274
$ws2.Range("C9").Value = 275.85
$ws2.Range("D9").Value = 272.95
$ws2.Range("E9").Value = 274.5
$ws2.Range("F9").Value = 2008463
$ws2.Range("G9").Value = 1397603
$ws2.Range("H9").Value = 0.4370769095372577
$ws2.Range("I9").Value = "PETRONET"
$ws2.Range("A10").Value = "PIIND"
$ws2.Range("B10").Value = 3460
$ws2.Range("C10").Value = 3483.8
$ws2.Range("D10").Value = 3430
$ws2.Range("E10").Value = 3440.1
$ws2.Range("F10").Value = 288764
$ws2.Range("G10").Value = 200324
$ws2.Range("H10").Value = 0.441484794632695
$ws2.Range("I10").Value = "PIIND"
$ws2.Range("A11").Value = "OIL"
$ws2.Range("B11").Value = 438
$ws2.Range("C11").Value = 438.9
$ws2.Range("D11").Value = 434.5
$ws2.Range("E11").Value = 436.2
$ws2.Range("F11").Value = 1092339
$ws2.Range("G11").Value = 760053
$ws2.Range("H11").Value = 0.4371879329467813
$ws2.Range("I11").Value = "OIL"
$ws2.Range("A12").Value = "COLPAL"
$ws2.Range("B12").Value = 2187.9
$ws2.Range("C12").Value = 2210
$ws2.Range("D12").Value = 2175.3
$ws2.Range("E12").Value = 2180
$ws2.Range("F12").Value = 387157
$ws2.Range("G12").Value = 260816
$ws2.Range("H12").Value = 0.4844066314950003
$ws2.Range("I12").Value = "COLPAL"
$ws2.Range("A13").Value = "ATGL"
$ws2.Range("B13").Value = 613.9
$ws2.Range("C13").Value = 621.4
$ws2.Range("D13").Value = 605.9
$ws2.Range("E13").Value = 608.05
$ws2.Range("F13").Value = 591539
$ws2.Range("G13").Value = 379035
$ws2.Range("H13").Value = 0.5606447953355231
$ws2.Range("I13").Value = "ATGL"
$ws2.Range("A14").Value = "POLYCAB"
$ws2.Range("B14").Value = 7700
$ws2.Range("C14").Value = 7740
$ws2.Range("D14").Value = 7630.5
$ws2.Range("E14").Value = 7649.5
$ws2.Range("F14").Value = 220950
$ws2.Range("G14").Value = 138316
$ws2.Range("H14").Value = 0.5974290754504179
$ws2.Range("I14").Value = "POLYCAB"
$ws2.Range("A15").Value = "KFINTECH"
$ws2.Range("B15").Value = 1083.1
$ws2.Range("C15").Value = 1098.9
$ws2.Range("D15").Value = 1077.8
$ws2.Range("E15").Value = 1080
$ws2.Range("F15").Value = 706819
$ws2.Range("G15").Value = 461199
$ws2.Range("H15").Value = 0.532568370703319
$ws2.Range("I15").Value = "KFINTECH"
